$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new row of data (2025-11-09) below the existing table.
# Column A holds the date as literal text (matching the existing rows which
# store dates as plain strings rather than Excel date serials), so a leading
# apostrophe is used to stop Excel from auto-converting it to a date value.
$ws.Range("A84").Value = "'11/09/2025"
$ws.Range("A84").Style = "Normal"
$ws.Range("B84").Value = 10059.82
